# "Fruta / hortaliza, semanal" — add this week's Femacal de La Calera /
# Repollo price report (two quality grades: Primera and Segunda) as two
# new rows inserted right before the existing "Crespo record" series,
# pushing the remaining history down by two rows (old A1:R316 -> A1:R318).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 288-289; Excel shifts rows 288..316 down to
# 290..318 and copies formatting (incl. the date style) from the row above.
$ws.Range("A288:A289").EntireRow.Insert()

# Row 288: new report, Calidad = Primera
$ws.Range("A288").Value = 3
$ws.Range("B288").Value = "Femacal de La Calera"
$ws.Range("C288").Value = "Coquimbo"
$ws.Range("D288").Value = 44449
$ws.Range("D288").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E288").Value = 5
$ws.Range("F288").Value = 100112006
$ws.Range("G288").Value = "Repollo"
$ws.Range("H288").Value = "Crespo record"
$ws.Range("I288").Value = "Primera"
$ws.Range("J288").Value = 3400
$ws.Range("K288").Value = 500
$ws.Range("L288").Value = 600
$ws.Range("M288").Value = 553
$ws.Range("N288").Value = "$/unidad"
$ws.Range("O288").Value = "Provincia de Quillota"
$ws.Range("P288").Value = 553
$ws.Range("Q288").Value = 1
$ws.Range("R288").Value = "Hortaliza"

# Row 289: new report, Calidad = Segunda
$ws.Range("A289").Value = 3
$ws.Range("B289").Value = "Femacal de La Calera"
$ws.Range("C289").Value = "Coquimbo"
$ws.Range("D289").Value = 44449
$ws.Range("D289").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E289").Value = 5
$ws.Range("F289").Value = 100112006
$ws.Range("G289").Value = "Repollo"
$ws.Range("H289").Value = "Crespo record"
$ws.Range("I289").Value = "Segunda"
$ws.Range("J289").Value = 1200
$ws.Range("K289").Value = 400
$ws.Range("L289").Value = 400
$ws.Range("M289").Value = 400
$ws.Range("N289").Value = "$/unidad"
$ws.Range("O289").Value = "Provincia de Quillota"
$ws.Range("P289").Value = 400
$ws.Range("Q289").Value = 1
$ws.Range("R289").Value = "Hortaliza"
